# WEAP/LEAP Time Step Control Test deck — apply the two reachable content
# edits from the target revision:
#   1. Slide 6 ("LEAP Model"): reposition the embedded picture.
#   2. Slide 7 ("Test Result: Unconstrained Source and Unconstrained Link
#      Test"): merge the "Unconstrained " / "Source " runs into a single
#      run "Unconstrained Source " (same run formatting, so this is a
#      pure text concatenation, not a visual change).

$p = $ppt.ActivePresentation

# --- 1. Slide 6 picture reposition -----------------------------------
$slide6 = $p.Slides.Item(6)
$pic = $slide6.Shapes.Item(4)

# Target EMU offsets: x=143626, y=1668150 (from x=102062, y=1596044).
# Shape.Left/.Top are in points (1 pt = 12700 EMU); a tiny epsilon nudges
# the float past the EMU rounding boundary so the stored value lands on
# the exact integer EMU target instead of one EMU short.
$pic.Left = (143626 / 12700.0) + 0.00004
$pic.Top  = (1668150 / 12700.0) + 0.00004

# --- 2. Slide 7 run merge --------------------------------------------
$slide7 = $p.Slides.Item(7)
$title = $slide7.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# "Test Result: " (13 chars) is followed by "Unconstrained " + "Source "
# (21 chars combined); rewriting that sub-range merges the two runs into
# one since they already share identical formatting.
$sub = $tr.Characters(14, 21)
$sub.Text = "Unconstrained Source "
